$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBD")
$r = $ws.Range("G11")
$v = $r.Value
Write-Host "value is: $v"
$v2 = $r.Value2
Write-Host "value2 is: $v2"
